$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Swap the contents of AB3:AI17 with AX3:BE17 (crosswalk cleanup of NLA location data)
$rangeLeft  = $ws.Range("AB3:AI17")
$rangeRight = $ws.Range("AX3:BE17")

$leftValues  = $rangeLeft.Value()
$rightValues = $rangeRight.Value()

$rangeLeft.Value  = $rightValues
$rangeRight.Value = $leftValues

# Update sheet view: keep the freeze at column B/C boundary (2 frozen columns),
# then move the selection to the new area of interest.
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C1").Select()
$win.FreezePanes = $true

$ws.Range("AX19:BF33").Select()
